$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07113107972035676
$ws.Range("H2").Value = -15.18660065395526
$ws.Range("I2").Value = -11.02553747096676
$ws.Range("G3").Value = 0.109712538079388
$ws.Range("H3").Value = -5.709373194861725
$ws.Range("G4").Value = -0.04334235139729783
$ws.Range("H4").Value = -180.9826133709635
$ws.Range("G5").Value = -0.02403659623448773
$ws.Range("H5").Value = 66.54692245025589
$ws.Range("G6").Value = 0.04351540542916248
$ws.Range("H6").Value = 22.42113169188307
$ws.Range("G7").Value = 0.05578991574085917
$ws.Range("H7").Value = 174.8935474988888
$ws.Range("G8").Value = -0.1441940543403089
$ws.Range("H8").Value = -2.194332918622467
$ws.Range("G9").Value = -0.1512036640870018
$ws.Range("H9").Value = -10.4587380833149
$ws.Range("G10").Value = -0.07488830920606399
$ws.Range("H10").Value = 30.67638290176361
$ws.Range("G11").Value = -0.09642597429516885
$ws.Range("H11").Value = -44.91242211587474
$ws.Range("G12").Value = -0.3995684129053694
$ws.Range("H12").Value = 3.763980947795051
$ws.Range("G13").Value = -0.4650545332056512
$ws.Range("H13").Value = -3.598834304604964
$ws.Range("G14").Value = -0.0664894990863929
$ws.Range("H14").Value = -31.05072820803433
$ws.Range("G15").Value = 0.04788052288441671
$ws.Range("H15").Value = 157.9026370629705
$ws.Range("G16").Value = 0.1425776452543689
$ws.Range("H16").Value = -2.359250223473684
$ws.Range("G17").Value = 0.1680684862885647
$ws.Range("H17").Value = 37.06283682605287
$ws.Range("G18").Value = 0.1348878766200803
$ws.Range("H18").Value = -2.577911819097585
$ws.Range("G19").Value = 0.1309136457753863
$ws.Range("H19").Value = 37.21634255255154
$ws.Range("G20").Value = 0.01784431859928216
$ws.Range("H20").Value = -30.3882254012952
$ws.Range("G21").Value = 0.05339173384659295
$ws.Range("H21").Value = -28.96928982374477
$ws.Range("G24").Value = 0.09143780897749496
$ws.Range("H24").Value = -8.959588765235587
$ws.Range("G25").Value = 0.1415073151844194
$ws.Range("H25").Value = -6.624512322224582
$ws.Range("G26").Value = 0.07970778036459751
$ws.Range("H26").Value = 0.7368985804903471
$ws.Range("G27").Value = 0.06401234847376806
$ws.Range("H27").Value = -35.89317157287365
$ws.Range("G28").Value = -0.2993869408769317
$ws.Range("H28").Value = -40.47864413226516
$ws.Range("G29").Value = -0.2379347487881293
$ws.Range("H29").Value = -15.89235849519372
$ws.Range("G30").Value = 0.06037406915809573
$ws.Range("H30").Value = 36.80449111647248
$ws.Range("G31").Value = 0.0347789738336989
$ws.Range("H31").Value = 32.07316335519585
$ws.Range("G32").Value = 0.1102841785658499
$ws.Range("H32").Value = 16.13688679775018
$ws.Range("G33").Value = 0.1192159242993803
$ws.Range("H33").Value = 14.66941175075153
$ws.Range("G34").Value = 0.05105453160520626
$ws.Range("H34").Value = 9.965691107010636
$ws.Range("G35").Value = 0.03704941282326985
$ws.Range("H35").Value = 389.0025490469271
$ws.Range("G36").Value = 0.06676471260712706
$ws.Range("H36").Value = 15.63743283426455
$ws.Range("G37").Value = 0.06387032189697642
$ws.Range("H37").Value = -9.179852179145161
$ws.Range("G38").Value = 0.007276797059882698
$ws.Range("H38").Value = -86.1093805248408
$ws.Range("G39").Value = 0.01960301020224836
$ws.Range("H39").Value = -5.469938322484507
$ws.Range("G40").Value = -0.02340155902684452
$ws.Range("H40").Value = -175.6046155691306
$ws.Range("G41").Value = 0.02113112751111922
$ws.Range("H41").Value = -40.23301627925598
$ws.Range("G42").Value = 0.1331046791015223
$ws.Range("H42").Value = -0.4349342014456794
$ws.Range("G43").Value = 0.1582357404591176
$ws.Range("H43").Value = 6.215353599696199
$ws.Range("G44").Value = 0.0008443675012952286
$ws.Range("H44").Value = 109.9203355372278
$ws.Range("G45").Value = -0.02633038165329555
$ws.Range("H45").Value = -139.8060587260944
$ws.Range("G46").Value = -0.02388432679770326
$ws.Range("H46").Value = -625.2958950413364
$ws.Range("G47").Value = -0.01695647487507642
$ws.Range("H47").Value = -82.75054614025316
$ws.Range("G48").Value = 0.05116786146639996
$ws.Range("H48").Value = 1.781801913812923
$ws.Range("G49").Value = 0.07136811459518631
$ws.Range("H49").Value = 8.027470416346187
$ws.Range("G50").Value = 0.1431706587327816
$ws.Range("H50").Value = -11.21622334485086
$ws.Range("G51").Value = 0.180493888884431
$ws.Range("H51").Value = 5.481069242374876
$ws.Range("G52").Value = -0.1621569569426365
$ws.Range("H52").Value = -1.075824066346579
$ws.Range("G53").Value = -0.1261687697314961
$ws.Range("H53").Value = -0.09065077166996861
$ws.Range("G54").Value = 0.1170795287136437
$ws.Range("H54").Value = 24.92161394800567
$ws.Range("G55").Value = 0.1260935950745711
$ws.Range("H55").Value = 11.50757208246073
$ws.Range("G56").Value = -0.02219236260563758
$ws.Range("H56").Value = -203.9728416699488
$ws.Range("G57").Value = -0.00007283743530589136
$ws.Range("H57").Value = 99.68147903476336
$ws.Range("G58").Value = 0.04274755747674795
$ws.Range("H58").Value = -24.19029176983394
$ws.Range("G59").Value = 0.0593586142394708
$ws.Range("H59").Value = -17.35077217410241
$ws.Range("G60").Value = 0.09035552577767574
$ws.Range("H60").Value = 29.1273700141456
$ws.Range("G61").Value = 0.08426843896438584
$ws.Range("H61").Value = 77.30725531779474
$ws.Range("G62").Value = 0.06649048740053821
$ws.Range("H62").Value = -8.867613855286608
$ws.Range("G63").Value = 0.06710569676811075
$ws.Range("H63").Value = 2.618290971747636
$ws.Range("G64").Value = -0.02849504783070192
$ws.Range("H64").Value = 31.19715412895823
$ws.Range("G65").Value = -0.004847950891731147
$ws.Range("H65").Value = 90.17115476768019
$ws.Range("G66").Value = 0.04488240263105538
$ws.Range("H66").Value = 137.0469046537676
$ws.Range("G67").Value = 0.04640519130130798
$ws.Range("H67").Value = 77.4584851328155
$ws.Range("G68").Value = -0.007597334327381687
$ws.Range("H68").Value = -1432.974703531967
$ws.Range("G69").Value = 0.001072836816825439
$ws.Range("H69").Value = 108.3086379069702
$ws.Range("G70").Value = -0.03070396398005286
$ws.Range("H70").Value = -11.86202194437141
$ws.Range("G71").Value = -0.02366936598007949
$ws.Range("H71").Value = 57.03661686538477
$ws.Range("G72").Value = -0.1385668657804331
$ws.Range("H72").Value = 6.584647331066837
$ws.Range("G73").Value = -0.1284915771111152
$ws.Range("H73").Value = 11.25649378570882
$ws.Range("G74").Value = 0.1323046001620441
$ws.Range("H74").Value = 4.961078678896199
$ws.Range("G75").Value = 0.1537748728016279
$ws.Range("H75").Value = 13.76239067385398
$ws.Range("G76").Value = -0.0495867365127257
$ws.Range("H76").Value = -43.9886626613908
$ws.Range("G77").Value = 0.005931709357930415
$ws.Range("H77").Value = 112.8428093426505
$ws.Range("G78").Value = 0.08921969374043801
$ws.Range("H78").Value = -3.202555251482371
$ws.Range("G79").Value = 0.09582699101861711
$ws.Range("H79").Value = -0.7019765390331975
$ws.Range("G80").Value = -0.1937122562041756
$ws.Range("H80").Value = -19.27150234380712
$ws.Range("G81").Value = -0.1587198902126401
$ws.Range("H81").Value = 26.67081499869976
$ws.Range("G82").Value = 0.1609783401555391
$ws.Range("H82").Value = 16.02353678638205
$ws.Range("G83").Value = 0.1955055024970945
$ws.Range("H83").Value = 18.76198730915315
$ws.Range("G84").Value = 0.03676295923534353
$ws.Range("H84").Value = 162.685865230001
$ws.Range("G85").Value = 0.06993785002077196
$ws.Range("H85").Value = 208.9239035030513
